# Regenerate merged AHB files
# 1) Rename header strings from *_old -> *_FV2310 and *_new -> *_FV2404
# 2) Turn the data range into an Excel Table (ListObject)
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2310Suffix = "_FV2310"
$fv2404Suffix = "_FV2404"

$headers = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $headers.Count; $i++) {
    # Columns A..J (1..10) carry the "_old" -> "_FV2310" headers
    $colOld = $i + 1
    $ws.Cells.Item(1, $colOld).Value = $headers[$i] + $fv2310Suffix

    # Columns L..U (12..21) carry the "_new" -> "_FV2404" headers (column K is "diff")
    $colNew = $i + 12
    $ws.Cells.Item(1, $colNew).Value = $headers[$i] + $fv2404Suffix
}

# Turn the whole used range (header + data, A1:U57) into an Excel Table named "Table1"
$tableRange = $ws.Range("A1:U57")
$listObject = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$listObject.Name = "Table1"

# Freeze the header row so it stays visible while scrolling
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
